$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "[-, -, -, 'MEC-3B-Coman. Hidraulicos']"

$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "['MCT-3A-Eletrohidráulica', -, -, -]"
$ws.Range("D3").Value = "['MEC-3B-Coman. Hidraulicos', -, -, -]"

$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "['MCT-3A-Eletrohidráulica', -, -, -]"
$ws.Range("D4").Value = "['MEC-3B-Coman. Hidraulicos', -, -, -]"

$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "['MCT-3A-Eletrohidráulica', -, -, -]"
$ws.Range("D6").Value = "-"

$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "['MCT-3A-Eletrohidráulica', -, -, -]"
$ws.Range("D7").Value = "-"

$ws.Range("C8").Value = "[-, -, 'MEC-3B-Coman. Hidraulicos', -]"

$ws.Range("E10").Value = "['MEC-3A-Cont.Lóg.Prog CLP', -, 'MEC-3A-Comandos Eletricos', -]"

$ws.Range("C11").Value = "MEC-2A-Máquinas Térmicas e de Fluxo"
$ws.Range("E11").Value = "-"

$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "['MEC-3A-Comandos Eletricos', -, 'MEC-3A-Cont.Lóg.Prog CLP', -]"
$ws.Range("E12").Value = "-"

$ws.Range("D14").Value = "['MEC-3A-Comandos Eletricos', -, 'MEC-3A-Cont.Lóg.Prog CLP', -]"
$ws.Range("E14").Value = "-"

$ws.Range("E15").Value = "-"

$ws.Range("D16").Value = "MEC-2A-Máquinas Térmicas e de Fluxo"
$ws.Range("E16").Value = "[-, 'MEC-3A-Cont.Lóg.Prog CLP', -, 'MEC-3A-Comandos Eletricos']"
